# Auto-generated edit script for kartu-piutang-2022.xlsx
# Fixes 'Taanggal' typo -> 'Tanggal', and populates the borrower
# header block (C5:C9 / F5:F9) plus the 12-row installment ledger
# (rows 12-23) for Asep Sukarsa.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Tanggal"
$ws.Range("C5").Value = 5080
$ws.Range("F5").Value = "lancar"
$ws.Range("C6").Value = "Asep Sukarsa"
$ws.Range("F6").Value = 9000000
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "2021-07-26"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2023-08-05"
$ws.Range("C8").Value = 3500000
$ws.Range("F8").Value = 75000000
$ws.Range("C9").Value = 84000000
$ws.Range("F9").Value = 9000000
$ws.Range("A12").Value = 1
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "2021-09-06"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "3,125,000"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "375,000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3,500,000"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "80,500,000"
$ws.Range("A13").Value = 2
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2021-10-07"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "3,125,000"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "375,000"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3,500,000"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "77,000,000"
$ws.Range("A14").Value = 3
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2021-11-21"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "3,125,000"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "375,000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3,500,000"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "73,500,000"
$ws.Range("A15").Value = 4
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "2021-12-16"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "3,125,000"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "375,000"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3,500,000"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "70,000,000"
$ws.Range("A16").Value = 5
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "2022-01-13"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "3,125,000"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "375,000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3,500,000"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "66,500,000"
$ws.Range("A17").Value = 6
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "2022-02-15"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "3,125,000"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "375,000"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3,500,000"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "63,000,000"
$ws.Range("A18").Value = 7
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2022-03-17"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "3,125,000"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "375,000"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3,500,000"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "59,500,000"
$ws.Range("A19").Value = 8
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "2022-04-24"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "3,125,000"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "375,000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3,500,000"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "56,000,000"
$ws.Range("A20").Value = 9
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "2022-05-31"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "3,125,000"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "375,000"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3,500,000"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "52,500,000"
$ws.Range("A21").Value = 10
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2022-06-30"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "3,125,000"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375,000"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3,500,000"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "49,000,000"
$ws.Range("A22").Value = 11
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "2022-07-15"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "3,125,000"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375,000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3,500,000"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "45,500,000"
$ws.Range("A23").Value = 12
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "2022-08-30"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "3,125,000"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "375,000"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3,500,000"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "42,000,000"
